$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 13:10"

# Update country rows: names (column A) and statistics (columns B-H)
# Country name changes reflect countries re-ranking in the sorted table;
# the B-H numbers are the refreshed case/death counts for whichever
# country now occupies that row.

$ws.Range("A13").Value = "India"
$rowVals = @(158959, 873, 67903, 86516, 0, 6, 4540)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(13, $i + 2).Value = $rowVals[$i]
}

$ws.Range("A59").Value = "Noruega"
$rowVals = @(8401, 0, 7727, 438, 0, 1, 236)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(59, $i + 2).Value = $rowVals[$i]
}

$ws.Range("A86").Value = "Bosnia y Herzegovina"
$rowVals = @(2462, 27, 1781, 528, 0, 2, 153)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(86, $i + 2).Value = $rowVals[$i]
}

$ws.Range("A100").Value = "Eslovenia"
$rowVals = @(1473, 2, 1356, 9, 0, 0, 108)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(100, $i + 2).Value = $rowVals[$i]
}

$ws.Range("A101").Value = "Kenia"
$rowVals = @(1471, 0, 408, 1008, 0, 0, 55)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(101, $i + 2).Value = $rowVals[$i]
}

$ws.Range("A122").Value = "Etiopia"
$rowVals = @(831, 100, 191, 633, 0, 1, 7)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(122, $i + 2).Value = $rowVals[$i]
}

$ws.Range("A123").Value = "Uruguay"
$rowVals = @(803, 0, 650, 131, 0, 0, 22)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(123, $i + 2).Value = $rowVals[$i]
}

$ws.Range("A124").Value = "Sierra Leona"
$rowVals = @(782, 0, 297, 440, 0, 0, 45)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(124, $i + 2).Value = $rowVals[$i]
}

$ws.Range("A125").Value = "Principado de Andorra"
$rowVals = @(763, 0, 676, 36, 0, 0, 51)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(125, $i + 2).Value = $rowVals[$i]
}

$ws.Range("A126").Value = "Nicaragua"
$rowVals = @(759, 0, 370, 354, 0, 0, 35)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(126, $i + 2).Value = $rowVals[$i]
}

$ws.Range("A127").Value = "Georgia"
$rowVals = @(738, 3, 573, 153, 0, 0, 12)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(127, $i + 2).Value = $rowVals[$i]
}

$ws.Range("A133").Value = "Madagascar"
$rowVals = @(656, 44, 154, 500, 0, 0, 2)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(133, $i + 2).Value = $rowVals[$i]
}

$ws.Range("A197").Value = "Curazao"
$rowVals = @(18, 0, 14, 3, 0, 0, 1)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(197, $i + 2).Value = $rowVals[$i]
}

$ws.Range("A198").Value = "Fiyi"
$rowVals = @(18, 0, 15, 3, 0, 0, 0)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(198, $i + 2).Value = $rowVals[$i]
}

$ws.Range("A199").Value = "Santa Lucia"
$rowVals = @(18, 0, 18, 0, 0, 0, 0)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(199, $i + 2).Value = $rowVals[$i]
}

$ws.Range("A201").Value = "Nueva Caledonia"
$rowVals = @(18, 0, 18, 0, 0, 0, 0)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(201, $i + 2).Value = $rowVals[$i]
}

$ws.Range("A210").Value = "Montserrat"
$rowVals = @(11, 0, 10, 0, 0, 0, 1)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(210, $i + 2).Value = $rowVals[$i]
}

$ws.Range("A211").Value = "Seychelles"
$rowVals = @(11, 0, 11, 0, 0, 0, 0)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(211, $i + 2).Value = $rowVals[$i]
}

$ws.Range("A213").Value = "Papua Nueva Guinea"
$rowVals = @(8, 0, 8, 0, 0, 0, 0)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(213, $i + 2).Value = $rowVals[$i]
}

$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$rowVals = @(8, 0, 7, 0, 0, 0, 1)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(214, $i + 2).Value = $rowVals[$i]
}

$ws.Range("A215").Value = "San Bartolome"
$rowVals = @(6, 0, 6, 0, 0, 0, 0)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(215, $i + 2).Value = $rowVals[$i]
}

$ws.Range("A216").Value = "Bonaire, San Eustaquio y Saba"
$rowVals = @(6, 0, 6, 0, 0, 0, 0)
for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item(216, $i + 2).Value = $rowVals[$i]
}
